# Update the "About" sheet and "Boundaries and methane sources" sheet
# with the new build/version string.

$wb = $excel.ActiveWorkbook

$oldText = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newText = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- "About" sheet updates ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value2 = "Version: " + $newText

$wsAbout.Range("A6").Value2 = "Recommended Citation:  " + '"Global Energy Monitor, Coal mine boundaries and methane sources for Broadmeadow Coal Mine, Australia, M0016, version ' + "'" + $newText + "'" + ". (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet updates ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 33; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # Column S is the 19th column
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}
